$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Rows 13-16: bump replicate count (M) 8 -> 9, fix dates in A14:A16
# ---------------------------------------------------------------
$ws.Range("M13").Value = 9
$ws.Range("M14").Value = 9
$ws.Range("M15").Value = 9
$ws.Range("M16").Value = 9

$ws.Range("A14").Value = 42578
$ws.Range("A15").Value = 42578
$ws.Range("A16").Value = 42578

# Row 16 totals: walltime per run changes, and the projected-total
# formula now multiplies by 4 (four clusters) instead of just summing.
$ws.Range("O16").Value = 0.625
$ws.Range("P16").Formula = "=SUM(N13:N16) *4"

# Q16 now displays as a full date/time AM-PM custom format
$ws.Range("Q16").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"

# ---------------------------------------------------------------
# New row 17: another simulation entry (6-core Macpro / Sim)
# ---------------------------------------------------------------
$ws.Range("A17").Value = 42578
$ws.Range("A17").NumberFormat = "d-mmm-yy"

$ws.Range("B17").Value = "6-core Macpro"
$ws.Range("B17").NumberFormat = "d-mmm-yy"

$ws.Range("C17").Value = "Sim"

$ws.Range("D17").Value = 25
$ws.Range("E17").Value = 1200
$ws.Range("F17").Value = 5000
$ws.Range("G17").Formula = "=E17*F17"

$ws.Range("H17").Formula = "=L17-K17"
$ws.Range("H17").NumberFormat = "h:mm"
$ws.Range("H17").HorizontalAlignment = -4152

$ws.Range("I17").Value = 1

$ws.Range("K17").Value = 0.61319444444444449
$ws.Range("K17").NumberFormat = "h:mm"

$ws.Range("L17").Value = 0.63194444444444442
$ws.Range("L17").NumberFormat = "h:mm"

$ws.Range("M17").Value = 9

$ws.Range("N17").Formula = "=M17*H17"
$ws.Range("N17").NumberFormat = "h:mm"

$ws.Range("O17").Value = 0.66666666666666696
$ws.Range("O17").NumberFormat = "h:mm"

$ws.Range("P17").Formula = "=H17"
$ws.Range("P17").NumberFormat = "h:mm"

$ws.Range("Q17").Formula = "=O17+P17"
$ws.Range("Q17").NumberFormat = "h:mm"

# ---------------------------------------------------------------
# View state: scrolled/selected a bit further right & down
# ---------------------------------------------------------------
$ws.Range("O21").Select()
